# "Journal de travail" worksheet update
#
# - Row 58 (D58): fix typo "fixs" -> "fixes" in the work-log comment.
# - Row 62 (C62/D62): log 8 hours of work and expand the "Frontend:
#   Correctifs" note with the extra details about config reset / error
#   handling; the cell needs to wrap onto two lines like the other
#   multi-line notes in this column.
# - Row 63 (A63/B63): add the next day's entry (date + "Implémentation").
# - Move the active selection from D65 to B65.
#
# (The "Temps [h]" total in C77 is a SUM formula over C2:C76, so it
# recalculates automatically once C62 is filled in.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Row 62 first, so the new shared string this introduces is appended
# before the one introduced by the D58 fix below (keeps string table
# ordering aligned with how Excel would have recorded the edits).
$ws.Range("D62").Value = "Frontend: Correctifs. Reset config + chargement de la configuration en db lors du démarrage de l'app. `nGestion des erreurs d'éxécution des modules (manager) et traitement des réponses (backend)"
$ws.Range("D62").WrapText = $true
$ws.Rows.Item(62).RowHeight = 34
$ws.Range("C62").Value = 8

# Fix the typo on row 58's comment.
$ws.Range("D58").Value = "Backend: Validation des request body, reset config default, fixes"

# New row for the following day.
$ws.Range("A63").Value = 45106
$ws.Range("B63").Value = "Implémentation"

# Move the active cell selection.
$ws.Range("B65").Select()
